$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new accelerometer sample was inserted as row 2 (shifting the former rows
# 2-21 down to rows 3-22), and 9 additional new samples were appended as rows
# 23-31. Rewrite every data cell directly with its final value so the
# existing cell formatting on rows 2-22 (no explicit style) is preserved and
# no incidental row-insert styles are introduced.
$ws.Range("A2").Value = 0.2003364562988281
$ws.Range("B2").Value = 0.6192827224731445
$ws.Range("C2").Value = -0.0097306966781616
$ws.Range("A3").Value = 0.1392936706542968
$ws.Range("B3").Value = 0.3245168924331665
$ws.Range("C3").Value = 0.0351060479879379
$ws.Range("A4").Value = 1.261228561401367
$ws.Range("B4").Value = 1.41911768913269
$ws.Range("C4").Value = -4.086977005004883
$ws.Range("A5").Value = 2.267774105072021
$ws.Range("B5").Value = -1.184367299079895
$ws.Range("C5").Value = -3.818997621536255
$ws.Range("A6").Value = -1.021368980407715
$ws.Range("B6").Value = 2.765882253646851
$ws.Range("C6").Value = 1.220625877380371
$ws.Range("A7").Value = -1.566243886947632
$ws.Range("B7").Value = 0.6303287744522095
$ws.Range("C7").Value = 0.5539150238037109
$ws.Range("A8").Value = -2.018197059631348
$ws.Range("B8").Value = 1.482012748718261
$ws.Range("C8").Value = 1.627924919128418
$ws.Range("A9").Value = -4.581077575683594
$ws.Range("B9").Value = 2.173830270767212
$ws.Range("C9").Value = 9.107954978942873
$ws.Range("A10").Value = -0.2751750946044922
$ws.Range("B10").Value = -0.9806771278381348
$ws.Range("C10").Value = -2.06553053855896
$ws.Range("A11").Value = 2.35319709777832
$ws.Range("B11").Value = 1.910999298095703
$ws.Range("C11").Value = -2.476747035980225
$ws.Range("A12").Value = -4.167366027832031
$ws.Range("B12").Value = 0.0987618193030357
$ws.Range("C12").Value = 3.953242778778076
$ws.Range("A13").Value = -1.787458419799805
$ws.Range("B13").Value = 1.655651211738586
$ws.Range("C13").Value = -5.035046577453613
$ws.Range("A14").Value = -9.838252067565918
$ws.Range("B14").Value = 3.984453201293945
$ws.Range("C14").Value = -6.098217010498047
$ws.Range("A15").Value = 6.411758422851562
$ws.Range("B15").Value = 1.583425164222717
$ws.Range("C15").Value = 7.352428436279297
$ws.Range("A16").Value = -2.261712551116944
$ws.Range("B16").Value = 0.8220813274383545
$ws.Range("C16").Value = 1.315514087677002
$ws.Range("A17").Value = -1.646389007568359
$ws.Range("B17").Value = 0.2190679311752319
$ws.Range("C17").Value = 0.9841623306274414
$ws.Range("A18").Value = -1.105591297149658
$ws.Range("B18").Value = 1.020219326019287
$ws.Range("C18").Value = 3.201179504394531
$ws.Range("A19").Value = -4.33466100692749
$ws.Range("B19").Value = -0.8289146423339844
$ws.Range("C19").Value = 6.12528133392334
$ws.Range("A20").Value = -1.558335304260254
$ws.Range("B20").Value = -0.159212052822113
$ws.Range("C20").Value = 1.605715155601502
$ws.Range("A21").Value = -0.9647946357727052
$ws.Range("B21").Value = 1.00678539276123
$ws.Range("C21").Value = -4.680802822113037
$ws.Range("A22").Value = -3.810809135437012
$ws.Range("B22").Value = 1.403007388114929
$ws.Range("C22").Value = 0.0495486259460449
$ws.Range("A23").Value = -1.585423946380615
$ws.Range("B23").Value = 2.060841083526612
$ws.Range("C23").Value = -2.507726192474365
$ws.Range("A24").Value = -5.486822128295898
$ws.Range("B24").Value = 2.457437515258789
$ws.Range("C24").Value = -1.076503276824951
$ws.Range("A25").Value = 3.813155174255371
$ws.Range("B25").Value = -5.157403945922852
$ws.Range("C25").Value = 7.194998264312744
$ws.Range("A26").Value = -3.507768154144287
$ws.Range("B26").Value = 2.501498937606812
$ws.Range("C26").Value = 0.7795240879058838
$ws.Range("A27").Value = 0.2215757369995117
$ws.Range("B27").Value = -0.4009582996368408
$ws.Range("C27").Value = 2.163901329040528
$ws.Range("A28").Value = 0.1625576019287109
$ws.Range("B28").Value = 1.34720504283905
$ws.Range("C28").Value = -0.6319388151168823
$ws.Range("A29").Value = 0.044438362121582
$ws.Range("B29").Value = -0.1398162841796875
$ws.Range("C29").Value = -0.8414495587348938
$ws.Range("A30").Value = -0.1983919143676757
$ws.Range("B30").Value = -0.413076639175415
$ws.Range("C30").Value = 0.2017757892608642
$ws.Range("A31").Value = 0.6235456466674805
$ws.Range("B31").Value = 1.087465167045593
$ws.Range("C31").Value = 0.6343502402305603
